$wb = $excel.ActiveWorkbook

$values = @{
    "2025" = @{
        A2 = 0
        B2 = 1037.265132737054
        C2 = 0
        D2 = 0
        E2 = 28926.05393052954
        F2 = 0
        G2 = 8095.925712661834
        H2 = 0
        I2 = 16171.06685703679
        J2 = 0
        K2 = 0
        L2 = 48492.22142001599
        M2 = 10595.37713982
        N2 = 7071.746316472206
        O2 = 6993.891045176661
    }
    "2030" = @{
        A2 = 0
        B2 = 4157.588990853394
        C2 = 0
        D2 = 0
        E2 = 45991.90904307188
        F2 = 0
        G2 = 8095.925712661834
        H2 = 0
        I2 = 37079.12819938764
        J2 = 0
        K2 = 0
        L2 = 54844.03303316472
        M2 = 17449.04999683176
        N2 = 9024.73482637065
        O2 = 9724.259001609127
    }
    "2035" = @{
        A2 = 2754.31755456332
        B2 = 6368.910634126893
        C2 = 0
        D2 = 0
        E2 = 57457.45307013817
        F2 = 0
        G2 = 8095.925712661834
        H2 = 0
        I2 = 52465.73681402855
        J2 = 0
        K2 = 0
        L2 = 54844.03303316472
        M2 = 21912.87293902603
        N2 = 13034.31228651121
        O2 = 12860.17252493772
    }
    "2040" = @{
        A2 = 2754.31755456332
        B2 = 6368.910634126893
        C2 = 0
        D2 = 0
        E2 = 57457.45307013817
        F2 = 0
        G2 = 8095.925712661834
        H2 = 0
        I2 = 52465.73681402855
        J2 = 0
        K2 = 0
        L2 = 54844.03303316472
        M2 = 21912.87293902603
        N2 = 13151.87171037628
        O2 = 12860.17252493772
    }
    "2045" = @{
        A2 = 5713.151062849596
        B2 = 6368.910634126893
        C2 = 0
        D2 = 0
        E2 = 57457.45307013817
        F2 = 0
        G2 = 8095.925712661834
        H2 = 0
        I2 = 52465.73681402855
        J2 = 0
        K2 = 0
        L2 = 54844.03303316472
        M2 = 21912.87293902603
        N2 = 13601.0893110962
        O2 = 14937.13200870449
    }
    "2050" = @{
        A2 = 5713.151062849596
        B2 = 6368.910634126893
        C2 = 0
        D2 = 0
        E2 = 57457.45307013817
        F2 = 0
        G2 = 8095.925712661834
        H2 = 0
        I2 = 52465.73681402855
        J2 = 0
        K2 = 0
        L2 = 54844.03303316472
        M2 = 21912.87293902603
        N2 = 13601.0893110962
        O2 = 14937.13200870449
    }
}

foreach ($sheetName in $values.Keys) {
    $sheetNameStr = [string]$sheetName
    $ws = $wb.Worksheets.Item($sheetNameStr)
    $cellValues = $values[$sheetName]
    foreach ($cellRef in $cellValues.Keys) {
        $cellRefStr = [string]$cellRef
        $ws.Range($cellRefStr).Value = $cellValues[$cellRef]
    }
}
